$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "In progress" backlog items were resolved (their follow/block work
# landed), so flip their Status column from "In progress" to "Done" -
# matching the formatting already used by the other "Done" rows.
$doneTemplate = $ws.Range("G32")
$doneTemplate.Copy()

$statusCells = @("G31", "G34", "G35", "G36", "G37")
foreach ($addr in $statusCells) {
    $cell = $ws.Range($addr)
    $cell.PasteSpecial(-4122)
    $cell.Value = "Done"
}

$excel.CutCopyMode = $false

# Restore the view to a plain, scrolled-to-top state with the new zoom
# level and selection the author left the workbook in.
$win = $excel.ActiveWindow
$win.Zoom = 115
$ws.Range("D8").Select()
